# Updated board from Excel
# Applies the tile-race board content updates to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# tile13 (row 15): was "Missing" / "Missing" / "missing.png"
$ws.Range("B15").Value = "Elder chaos"
$ws.Range("C15").Value = "Obtain a pieve of elder chaos robes"
$ws.Range("D15").Value = "Elder_chaos_top.png"

# tile18 (row 20): was "Roll again" / "Roll again" / "dice.png"
$ws.Range("B20").Value = "Blood shard"
$ws.Range("C20").Value = "Obtain a blood shard"
$ws.Range("D20").Value = "Blood_shard.png"

# tile32 (row 34): must-hit flag flipped to FALSE
$ws.Range("E34").Value = $false

# tile34 (row 36): wrath talisman count 3 -> 5
$ws.Range("C36").Value = "Obtain 5 wrath talismans from Vorkath"

# tile36 (row 38): reworded description + must-hit flag flipped to TRUE
$ws.Range("C38").Value = "Complete a deathless CM raid (Immortal raid team) and 350 ToA with no deaths (Something of an expert myself)"
$ws.Range("E38").Value = $true

# tile37 (row 39): was "Missing" / "Missing" / "missing.png"
$ws.Range("B39").Value = "Champion scroll"
$ws.Range("C39").Value = "Obtain any champion scroll"
$ws.Range("D39").Value = "champion_scroll.png"

# tile38 (row 40): was "Missing" / "Missing" / "missing.png"
$ws.Range("B40").Value = "Chaos ely"
$ws.Range("C40").Value = "Get pet or dragon pickaxe from chaos elemental"
$ws.Range("D40").Value = "Dragon_pick.png"

# tile41 (row 43): was "Roll again" / "Roll again" / "dice.png"
$ws.Range("B43").Value = "Common GWD"
$ws.Range("C43").Value = "Get any of: Zamorakian Spear, Steam Battlestaff, Saradomin Sword, or Saradomin's light"
$ws.Range("D43").Value = "Zamorakian_spear.png"

# tile43 (row 45): was "Common GWD" set, now becomes "Zulrah fang" set
$ws.Range("B45").Value = "Zulrah fang"
$ws.Range("C45").Value = "Get a fang item from Zulrah (Magic Fang or Tanzanite fang)"
$ws.Range("D45").Value = "magic_fang.png"

# tile47 (row 49): was "Zulrah fang" set, now becomes "GWD armor" set
$ws.Range("B49").Value = "GWD armor"
$ws.Range("C49").Value = "Obtain a armor piece from GWD (bandos boots, chest and plate + arma chest, legs and helm)"
$ws.Range("D49").Value = "Arma_chest.png"
